$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "dfghj"
$ws.Range("C1").Value = "qsdfbn"
$ws.Range("D1").Value = "asdfg"
$ws.Range("E1").Value = "WDFGH"
$ws.Range("F2").Value = "dfghjk"
$ws.Range("G2").Value = "ghjk"
$ws.Range("H2").Value = "sdfjkl"
$ws.Range("I2").Value = "dfghjkl"

$ws.Range("I2").Select()
